$d = $word.ActiveDocument

# 1) Remove the first paragraph entirely:
#    "Możliwość wykupu akcji przez spsółkę"
$d.Paragraphs.Item(1).Range.Delete()

# 2) Remove the paragraph "Poprawa tabel w okienkach z informacjami"
#    (after step 1 this is now paragraph index 2)
$d.Paragraphs.Item(2).Range.Delete()

# 3) Replace the paragraph "Wywalenie zbędnych metod z interfejsu allinstances"
#    (now paragraph index 2) with the new content:
#    "Usunięcie metod display, setvalues, getOutputString"
$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Usunięcie metod display, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>setvalues</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getOutputString</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$d.Paragraphs.Item(2).Range.InsertXML($newParaXml)

# 4) Change the text of the last paragraph from "Wątek surowców i walut" to
#    "Wielkość okna startowego do ogarnięcia", keeping the _GoBack bookmark intact.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$lastPara.Range.Find.Execute("Wątek surowców i walut", $false, $false, $false, $false, $false, $true, 1, $false, "Wielkość okna startowego do ogarnięcia", 2)
